$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-10 -> 2023-09-11, i.e. 45179 -> 45180) for every data row
# (rows 2 through 171).
for ($row = 2; $row -le 171; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current + 1
    }
}
